$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell contents for rows 13-25 to match the new Ementa layout ---

$ws.Range("A13").Clear()

$ws.Range("B13").Value = '8767640 - Eduardo Ferro dos Santos'

$ws.Range("C13").Value = '8767640 - Eduardo Ferro dos Santos'

$ws.Range("A14").Value = 'Programa resumido:'

$ws.Range("B3").Copy($ws.Range("B14"))
$ws.Range("B14").Value = 'Controle e automação; Robótica; Domótica; Sistemas Supervisórios, Pneumática, Hidráulica, CLP'

$ws.Range("C3").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 'Controle e automação; Robótica; Domótica; Sistemas Supervisórios, Pneumática, Hidráulica, CLP'

$ws.Range("A15").Value = 'Short syllabus:'

$ws.Range("B15").Clear()

$ws.Range("C15").Clear()

$ws.Range("A16").Value = 'Programa:'

$ws.Range("B3").Copy($ws.Range("B16"))
$ws.Range("B16").Value = 'Introdução aos princípios de controle e automação; Fundamentos da Robótica; Fundamentos da Domótica;  Introdução a Sistemas Supervisórios, Princípios da Automação Pneumática, Hidráulica, Introdução aos Controladores Lógicos Programáveis.'

$ws.Range("C3").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 'Introdução aos princípios de controle e automação; Fundamentos da Robótica; Fundamentos da Domótica;  Introdução a Sistemas Supervisórios, Princípios da Automação Pneumática, Hidráulica, Introdução aos Controladores Lógicos Programáveis.'

$ws.Range("A17").Value = 'Syllabus:'

$ws.Range("A18").Value = 'Avaliação:'

$ws.Range("B18").Clear()

$ws.Range("C18").Clear()

$ws.Range("A19").Value = 'Método:'

$ws.Range("A20").Value = 'Critério:'

$ws.Range("A21").Value = 'Norma de recuperação:'

$ws.Range("A22").Value = 'Bibliografia:'

$ws.Range("B3").Copy($ws.Range("B22"))
$ws.Range("B22").Value = 'Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) NISE, N. S., “Engenharia de Sistemas de Controle”, 3ª ed., LTC, 2002. OGATA, K., “Engenharia de Controle Moderno”, 4ª ed., Prentice-Hall do Brasil, 2003. Tutoriais disponibilizados pelo professor BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U. B.. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p. CAPELLI, A. Automação Industrial: controle de movimento e processos contínuos. São Paulo: Érica, 2006. SILVEIRA, P. R. da; SANTOS, W. E. Automação e controle discreto. 3. ed. São Paulo: Érica, 1998. MORAES, C. C.; CATRUCCI, P. Engenharia de automação industrial. 2. ed. Rio de Janeiro: LTC, 2007. GIORGINI, M. Automação aplicada: descrição e implementação de sistemas sequencias com PLC''s. 5. ed. São Paulo: Érica, 2003.'

$ws.Range("C3").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 'Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) NISE, N. S., “Engenharia de Sistemas de Controle”, 3ª ed., LTC, 2002. OGATA, K., “Engenharia de Controle Moderno”, 4ª ed., Prentice-Hall do Brasil, 2003. Tutoriais disponibilizados pelo professor BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U. B.. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p. CAPELLI, A. Automação Industrial: controle de movimento e processos contínuos. São Paulo: Érica, 2006. SILVEIRA, P. R. da; SANTOS, W. E. Automação e controle discreto. 3. ed. São Paulo: Érica, 1998. MORAES, C. C.; CATRUCCI, P. Engenharia de automação industrial. 2. ed. Rio de Janeiro: LTC, 2007. GIORGINI, M. Automação aplicada: descrição e implementação de sistemas sequencias com PLC''s. 5. ed. São Paulo: Érica, 2003.'

$ws.Range("A3").Copy($ws.Range("A23"))
$ws.Range("A23").Value = 'Requisitos:'

$ws.Range("B23").Clear()

$ws.Range("C23").Clear()

$ws.Range("B24").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"

$ws.Range("C24").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"

$ws.Range("B3").Copy($ws.Range("B25"))
$ws.Range("B25").Value = "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)`n"

$ws.Range("C3").Copy($ws.Range("C25"))
$ws.Range("C25").Value = "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)`n"

# --- Fix up row heights to match the new layout ---

$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(25).RowHeight = 30
